$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- ATOC (column G) updates ---
# Week 1 (G3): rename video file to include .mp4
$ws.Range("G3").Value = "Topic: name- Week 1,lectures- 1 Lecture, duration- 01:20;`nVideo: link- https://drive.google.com/file/d/19dipmMmEiCQmlntPtJAtKF_DfoD_sufk/preview, name- ATC Week#1.mp4, duration- 01:20;`nSlides: slide- Course.pdf;`nSlides: slide- Lec1.pdf;`nTopics: Finite state machine, Applications of Finite Automata, Proving techniques"

# Week 2 (G4): rename video file to include .mp4
$ws.Range("G4").Value = "Topic: name- Week 2,lectures- 1 Lecture, duration- 01:44;`nVideo: link- https://drive.google.com/file/d/1LLPqYvPSZruiDnWQfo4C6cxPAwNAJBAc/preview, name- ATC Week#2.mp4, duration- 01:44;`nSlides: slide- Lec1.pdf;"

# Week 3 (G5): split lecture names into Part1/Part2
$ws.Range("G5").Value = "Topic: name- Week 3,lectures- 2 Lectures, duration- 01:49;`nVideo: link- https://drive.google.com/file/d/15puMyWTcN0p2nEstcIl5CfRdY9_QJMQ9/preview, name- ATC Week#3Part1.mp4, duration- 00:31;`nVideo: link- https://drive.google.com/file/d/1LHMHVg2vCtKk4WZw6evB8Ly6JOxUGj8X/preview, name- ATC Week#3Part2.mp4, duration- 01:18;`nSlides: slide- Lec1.pdf;"

# Week 4 (G6): brand new entry, copy formatting from G5 first
$ws.Range("G5").Copy()
$ws.Range("G6").PasteSpecial(-4122)
$ws.Range("G6").Value = "Topic: name- Week 4,lectures- 3 Lectures, duration- 01:40;`nVideo: link- https://drive.google.com/file/d/1dWC0HSI8lJ4ASxDgWuTYJ3bl1GrZ-Z-5/preview, name- ATOC Week#4Part1.mp4, duration- 00:50;`nVideo: link- https://drive.google.com/file/d/1qLxWUU9e9HeNY8DNLhalPa_CQbaxo8k_/preview, name- ATOC Week#4Part2.mp4, duration- 00:48;`nVideo: link- https://drive.google.com/file/d/1ZKP1-pdg6o_ifZ4GOkHRCjUUrnaqArza/preview, name- ATOC Week#4Part3.mp4, duration- 00:02;"

# --- AOS (column H) updates ---
# Week 2 (H4): rename video file
$ws.Range("H4").Value = "Topic: name- Week 2,lectures- 1 Lecture, duration- 01:33;`nVideo: link- https://drive.google.com/file/d/18qBlYxrN38HZh1pjIysBGtg8SF1JSLY0/preview, name- AOS Week#2part1.mp4, duration- 01:33;`nSlides: slide- Chapter 3.ppt;"

# Week 3 (H5): brand new entry, copy formatting from H4 first
$ws.Range("H4").Copy()
$ws.Range("H5").PasteSpecial(-4122)
$ws.Range("H5").Value = "Topic: name- Week 3,lectures- 2 Lectures, duration- 01:38;`nVideo: link- https://drive.google.com/file/d/1CEA0dNMZFnAwwcInBt37frFSakNJSyuO/preview, name- AOS Week#3part1.mp4, duration- 01:06;`nVideo: link- https://drive.google.com/file/d/19aBooNe8gvSiXqm-yxEMN2AHg0KthhYP/preview, name- AOS Week#3part2.mp4, duration- 00:32;"

# --- sheet view: scroll / selection moved to F6/G6 ---
$ws.Application.ActiveWindow.ScrollRow = 6
$ws.Application.ActiveWindow.ScrollColumn = 6
$ws.Range("G6").Select()
